$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert three new rows above row 4 (old rows 4-22 shift down to 7-25).
#    This also pulls the SUM formulas in what becomes row 25 along for the
#    ride (Excel auto-expands C2:C21 / D2:D21 to C2:C24 / D2:D24).
# ---------------------------------------------------------------------------
$ws.Rows("4:6").Insert()

# Copy the formatting (styles) of row 3 (which already carries the date
# number-format on E/F) down into the three freshly inserted rows so they
# pick up the same per-column styles used throughout the table.
$ws.Range("A3:I3").Copy()
$ws.Range("A4:I6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row heights for the new rows (from the authored change).
$ws.Rows("4").RowHeight = 31.5
$ws.Rows("5").RowHeight = 30
$ws.Rows("6").RowHeight = 30

# ---------------------------------------------------------------------------
# 2. Populate the three new task rows (S.No 3, 4, 5 - the Jawbone Login
#    workflow tasks added by this commit).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Make lib for QRScan and iBeacons"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 24
$ws.Range("E4").Value = "12/11/2014"
$ws.Range("F4").Value = "12/16/2014"
$ws.Range("G4").Value = "liubin&zhangxiaodong"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "Completed"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Investigat the Jawbone & Nick API, select wearable device"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 32
$ws.Range("E5").Value = "12/17/2014"
$ws.Range("F5").Value = "12/23/2014"
$ws.Range("G5").Value = "liubin&zhangxiaodong"
$ws.Range("H5").Value = "IDC want to use Nick+, but SDK can't get"
$ws.Range("I5").Value = "Completed"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Create Jawbone Login View base on app of Q1 and Jawbone developer account"
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 48
$ws.Range("E6").Value = "12/24/2014"
$ws.Range("F6").Value = "12/31/2014"
$ws.Range("G6").Value = "liubin"
$ws.Range("H6").Value = "Send the Screenshot to IDC which is Jawbone Login View."
$ws.Range("I6").Value = "Completed"

# ---------------------------------------------------------------------------
# 3. Updates to the pre-existing rows, now shifted down by 3 (rows 7-14).
#    The row insert carried the old S.No (3,4,5,6,7,8,9,10) along with the
#    row content, but the task-sequence column (A) needs renumbering to
#    account for the 3 new tasks inserted earlier in the list.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13

# Row 7 - "crossplatform(iOS & Android) design..." (was row 4)
$ws.Range("G7").Value = "liubin&zhangxiaodong"
$ws.Range("E3:F3").Copy()
$ws.Range("E7:F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E7").Value = "01/05/2015"
$ws.Range("F7").Value = "to the end"

# Row 8 - "design Hybris data structure..." (was row 5)
$ws.Range("E3:F3").Copy()
$ws.Range("E8:F8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E8").Value = "?"
$ws.Range("F8").Value = "?"

# Row 9 - " create new views in native iOS app with UI design" (was row 6)
$ws.Range("G9").Value = "liubin&zhangxiaodong"
$ws.Range("E3:F3").Copy()
$ws.Range("E9:F9").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E9").Value = "01/05/2015"
$ws.Range("F9").Value = "to the end"

# Row 10 - "Based on the app of Q1..." (was row 7): effort cleared out
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()

# Row 11 - "improve the iOS app, collecting..." (was row 8): effort cleared, owner updated
$ws.Range("C11").ClearContents()
$ws.Range("D11").ClearContents()
$ws.Range("G11").Value = "liubin&zhangxiaodong"

# Row 12 - "improve the iOS app, coworking..." (was row 9): effort cleared, owner updated
$ws.Range("C12").ClearContents()
$ws.Range("D12").ClearContents()
$ws.Range("G12").Value = "liubin&zhangxiaodong"

# Row 13 - "integration with Hybris and local test" (was row 10): effort cleared, owner updated
$ws.Range("C13").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("G13").Value = "liubin&zhangxiaodong"

# Row 14 - "demo and test with IDC team" (was row 11): effort cleared
$ws.Range("C14").ClearContents()
$ws.Range("D14").ClearContents()

# ---------------------------------------------------------------------------
# 4. Sheet view selection moves to E10 (matches authored change).
# ---------------------------------------------------------------------------
$ws.Range("E10").Select()
